$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("G6").Value = 2.2
$ws.Range("I6").Value = 3.25

# Row 7
$ws.Range("O7").Value = 1.22
$ws.Range("P7").Value = 4.33
$ws.Range("Q7").Value = 1.7
$ws.Range("R7").Value = 2.1

# Row 10
$ws.Range("W10").Value = 8.5
$ws.Range("AE10").Value = 17
$ws.Range("AG10").Value = 6.5
$ws.Range("AK10").Value = 21
$ws.Range("AM10").Value = 401

# Row 11
$ws.Range("M11").Value = 1.11
$ws.Range("N11").Value = 6.5
